$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.066.09"
$ws.Range("E2").Value = "  +6.11%  "
$ws.Range("D3").Value = "3.681.41"
$ws.Range("E3").Value = "  +18.41%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'599.09"
$ws.Range("E5").Value = "  +3.52%  "
$ws.Range("D6").Value = "'184.19"
$ws.Range("E6").Value = "  +5.98%  "
$ws.Range("D7").Value = "3.678.47"
$ws.Range("E7").Value = "  +18.40%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  +4.15%  "
$ws.Range("E10").Value = "  +7.01%  "
$ws.Range("D11").Value = "'6.63"
$ws.Range("E11").Value = "  +3.78%  "
$ws.Range("E12").Value = "  +5.02%  "
$ws.Range("D13").Value = "'40.17"
$ws.Range("E13").Value = "  +11.53%  "
$ws.Range("E14").Value = "  +5.48%  "
$ws.Range("D15").Value = "4.299.04"
$ws.Range("E15").Value = "  +18.80%  "
$ws.Range("D16").Value = "71.109.29"
$ws.Range("E16").Value = "  +6.36%  "
$ws.Range("D17").Value = "3.679.39"
$ws.Range("E17").Value = "  +18.81%  "
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "'7.50"
$ws.Range("E19").Value = "  +7.02%  "
$ws.Range("D20").Value = "'16.96"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "'514.53"
$ws.Range("E21").Value = "  +6.00%  "
$ws.Range("D22").Value = "'9.18"
$ws.Range("E22").Value = "  +17.19%  "
$ws.Range("D23").Value = "'0.744"
$ws.Range("E23").Value = "  +7.38%  "
$ws.Range("D24").Value = "'87.61"
$ws.Range("E24").Value = "  +4.83%  "
$ws.Range("D25").Value = "'2.50"
$ws.Range("E25").Value = "  +11.70%  "
$ws.Range("D26").Value = "'13.54"
$ws.Range("E26").Value = "  +5.88%  "
$ws.Range("D27").Value = "'11.05"
$ws.Range("E27").Value = "  +8.68%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +10.99%  "
$ws.Range("D30").Value = "'8.20"
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("E31").Value = "  +18.25%  "
$ws.Range("E32").Value = "  +6.93%  "
$ws.Range("D33").Value = "'31.59"
$ws.Range("E33").Value = "  +12.33%  "
$ws.Range("D34").Value = "'0.116"
$ws.Range("E34").Value = "  +3.59%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  +8.89%  "
$ws.Range("E37").Value = "  +7.28%  "
$ws.Range("D38").Value = "'0.347"
$ws.Range("E38").Value = "  +11.51%  "
$ws.Range("E39").Value = "  +9.84%  "
$ws.Range("D40").Value = "'51.12"
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("E41").Value = "  +3.79%  "
$ws.Range("D42").Value = "'45.25"
$ws.Range("E42").Value = "  -6.25%  "
$ws.Range("D43").Value = "3.141.89"
$ws.Range("E43").Value = "  +12.13%  "
$ws.Range("E44").Value = "  +6.50%  "
$ws.Range("D45").Value = "'420.61"
$ws.Range("E45").Value = "  +13.16%  "
$ws.Range("D46").Value = "'2.78"
$ws.Range("E46").Value = "  +4.50%  "
$ws.Range("D47").Value = "'0.0369"
$ws.Range("E47").Value = "  +6.39%  "
$ws.Range("D48").Value = "'28.42"
$ws.Range("E48").Value = "  +15.79%  "
$ws.Range("D49").Value = "'137.45"
$ws.Range("E49").Value = "  +2.52%  "
$ws.Range("E51").Value = "  +12.18%  "